# "DT added label option"
#
# The backwardElimination workbook holds, on each of its 28 sheets, a B2
# cell with the verbatim text dump of a statsmodels OLS summary (one
# elimination step per sheet). The summary was regenerated - same fitted
# numbers, but a new run timestamp - so the "Date:" and "Time:" header
# lines changed on every sheet while everything else (coefficients,
# R-squared, AIC/BIC, etc.) stayed identical.
#
# New timestamp observed in the refreshed run: Wed, 08 Jan 2020, with the
# wall clock ticking from 19:07:36 to 19:07:37 partway through the 28
# steps (exactly like the original run ticked from 21:22:31 to 21:22:32).

$wb = $excel.ActiveWorkbook

$newDate = "Wed, 08 Jan 2020"

# Per-sheet (elimination step) "Time:" value, in sheet order (same order
# as Worksheets.Item(1..28) / sheet1.xml..sheet28.xml).
$times = @(
    "19:07:36", "19:07:36", "19:07:36", "19:07:36", "19:07:36",
    "19:07:36", "19:07:36", "19:07:36", "19:07:36", "19:07:36",
    "19:07:36", "19:07:36", "19:07:36", "19:07:36", "19:07:36",
    "19:07:36", "19:07:36", "19:07:37", "19:07:37", "19:07:37",
    "19:07:37", "19:07:37", "19:07:37", "19:07:37", "19:07:37",
    "19:07:37", "19:07:37", "19:07:37"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $text = $ws.Range("B2").Value()
    if ($text -eq $null) { continue }

    $newTime = $times[$i - 1]

    # Keep column alignment identical - only the date/time token itself
    # changes, the fixed-width padding around "Date:"/"Time:" is untouched.
    $text = $text -replace "(Date:\s+)\w+, \d\d \w+ \d\d\d\d", ('${1}' + $newDate)
    $text = $text -replace "(Time:\s+)\d\d:\d\d:\d\d", ('${1}' + $newTime)

    $ws.Range("B2").Value = $text

    # Writing the cell makes the host recompute row 2's height; pin it
    # back to the sheet's original (Excel's own max) row height so the
    # row geometry doesn't drift from the source file.
    $ws.Rows.Item(2).RowHeight = 409.5
}
